$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 64668.668
$ws.Range("I13").Value = 39000
$ws.Range("J13").Value = 77503
$ws.Range("K13").Value = 39000
$ws.Range("L13").Value = 77503
$ws.Range("M13").Value = -38831
$ws.Range("N13").Value = -77841

# Row 63
$ws.Range("H63").Value = 39263
$ws.Range("J63").Value = 39263
$ws.Range("L63").Value = 39263
$ws.Range("N63").Value = -40511

# Row 66
$ws.Range("H66").Value = 39263
$ws.Range("J66").Value = 39263
$ws.Range("L66").Value = 117789
$ws.Range("N66").Value = -124029

# Row 121
$ws.Range("H121").Value = 960
$ws.Range("J121").Value = 1216
$ws.Range("L121").Value = 3648
$ws.Range("N121").Value = -7142

# Row 131
$ws.Range("H131").Value = 1574.591
$ws.Range("I131").Value = 1241.1666
$ws.Range("J131").Value = 3075
$ws.Range("K131").Value = 3723.4998
$ws.Range("L131").Value = 9225
$ws.Range("M131").Value = 1316.5002
$ws.Range("N131").Value = -19305

# Row 132
$ws.Range("H132").Value = 3863602.2
$ws.Range("I132").Value = 5717357
$ws.Range("J132").Value = 1613.4166
$ws.Range("K132").Value = 17152071
$ws.Range("L132").Value = 4840.2498
$ws.Range("M132").Value = -17149541
$ws.Range("N132").Value = -9900.2498

# Row 135
$ws.Range("H135").Value = 842.81396
$ws.Range("I135").Value = 593.525
$ws.Range("K135").Value = 5341.724999999999
$ws.Range("M135").Value = -2806.724999999999

# Row 137
$ws.Range("H137").Value = 2802.818
$ws.Range("I137").Value = 3398.3333
$ws.Range("K137").Value = 10194.9999
$ws.Range("M137").Value = -7644.999899999999

# Row 141
$ws.Range("H141").Value = 1728.8182
$ws.Range("I141").Value = 994.79486
$ws.Range("K141").Value = 2984.38458
$ws.Range("M141").Value = 2195.61542

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1090.13
$ws.Range("I32").Value = 1010.449
$ws.Range("J32").Value = 4994.5
$ws.Range("K32").Value = 1010.449
$ws.Range("L32").Value = 4994.5
$ws.Range("M32").Value = -723.449
$ws.Range("N32").Value = -5568.5

# Row 61
$ws.Range("H61").Value = 1367.4872
$ws.Range("I61").Value = 854.125
$ws.Range("J61").Value = 3714.2856
$ws.Range("K61").Value = 854.125
$ws.Range("L61").Value = 3714.2856
$ws.Range("M61").Value = -642.125
$ws.Range("N61").Value = -4138.2856

# Row 74
$ws.Range("H74").Value = 1080.4231
$ws.Range("I74").Value = 926.86365
$ws.Range("J74").Value = 1925
$ws.Range("K74").Value = 926.86365
$ws.Range("L74").Value = 1925
$ws.Range("M74").Value = -52.86365000000001
$ws.Range("N74").Value = -3673

# Row 77
$ws.Range("H77").Value = 1080.4231
$ws.Range("I77").Value = 926.86365
$ws.Range("J77").Value = 1925
$ws.Range("K77").Value = 4634.31825
$ws.Range("L77").Value = 9625
$ws.Range("M77").Value = -266.3182500000003
$ws.Range("N77").Value = -18361

# Row 136
$ws.Range("H136").Value = 1367.4872
$ws.Range("I136").Value = 854.125
$ws.Range("J136").Value = 3714.2856
$ws.Range("K136").Value = 2562.375
$ws.Range("L136").Value = 11142.8568
$ws.Range("M136").Value = -12.375
$ws.Range("N136").Value = -16242.8568

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1078.091
$ws.Range("I94").Value = 711.05884
$ws.Range("J94").Value = 2326
$ws.Range("K94").Value = 711.05884
$ws.Range("L94").Value = 2326
$ws.Range("M94").Value = -260.05884
$ws.Range("N94").Value = -3228

# Row 105
$ws.Range("H105").Value = 3379.8572
$ws.Range("I105").Value = 2846.0908
$ws.Range("J105").Value = 5337
$ws.Range("K105").Value = 2846.0908
$ws.Range("L105").Value = 5337
$ws.Range("M105").Value = -1099.0908
$ws.Range("N105").Value = -8831

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8775421
$ws.Range("I31").Value = 2581.3333
$ws.Range("J31").Value = 23814574
$ws.Range("K31").Value = 2581.3333
$ws.Range("L31").Value = 23814574
$ws.Range("M31").Value = -2286.3333
$ws.Range("N31").Value = -23815164

# Row 34
$ws.Range("H34").Value = 8775421
$ws.Range("I34").Value = 2581.3333
$ws.Range("J34").Value = 23814574
$ws.Range("K34").Value = 2581.3333
$ws.Range("L34").Value = 23814574
$ws.Range("M34").Value = -2379.3333
$ws.Range("N34").Value = -23814978

$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 39.75
$ws.Range("I10").Value = 39.75
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 119.25
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 19.75
$ws.Range("N10").ClearContents()

# Row 11
$ws.Range("H11").Value = 5883323
$ws.Range("I11").Value = 142.14285
$ws.Range("K11").Value = 426.42855
$ws.Range("M11").Value = -286.42855

# Row 13
$ws.Range("H13").Value = 800.1429000000001
$ws.Range("I13").Value = 700.3333
$ws.Range("J13").Value = 875
$ws.Range("K13").Value = 2100.9999
$ws.Range("L13").Value = 2625
$ws.Range("M13").Value = -1932.9999
$ws.Range("N13").Value = -2961

# Row 15
$ws.Range("H15").Value = 500
$ws.Range("I15").Value = 500
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 1500
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = -1360
$ws.Range("N15").Value = -1780

# Row 131
$ws.Range("H131").Value = 326790.38
$ws.Range("J131").Value = 397767.4
$ws.Range("L131").Value = 1193302.2
$ws.Range("N131").Value = -1203382.2

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 14707464
$ws.Range("I113").Value = 25001800
$ws.Range("J113").Value = 1271.1428
$ws.Range("K113").Value = 25001800
$ws.Range("L113").Value = 1271.1428
$ws.Range("M113").Value = -24999630
$ws.Range("N113").Value = -5611.1428

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1825.6
$ws.Range("I7").Value = 1277.8462
$ws.Range("J7").Value = 2842.8572
$ws.Range("K7").Value = 1277.8462
$ws.Range("L7").Value = 2842.8572
$ws.Range("M7").Value = -1165.8462
$ws.Range("N7").Value = -3066.8572

# Row 126
$ws.Range("H126").Value = 1825.6
$ws.Range("I126").Value = 1277.8462
$ws.Range("J126").Value = 2842.8572
$ws.Range("K126").Value = 3833.5386
$ws.Range("L126").Value = 8528.571599999999
$ws.Range("M126").Value = -1363.5386
$ws.Range("N126").Value = -13468.5716

# Row 127
$ws.Range("H127").Value = 35300
$ws.Range("J127").Value = 35300
$ws.Range("L127").Value = 35300
$ws.Range("N127").Value = -45220

# Row 133
$ws.Range("H133").Value = 20902.6
$ws.Range("J133").Value = 20902.6
$ws.Range("L133").Value = 20902.6
$ws.Range("N133").Value = -25962.6

# Row 136
$ws.Range("H136").Value = 3512.9092
$ws.Range("I136").Value = 3814.0303
$ws.Range("J136").Value = 2609.5454
$ws.Range("K136").Value = 11442.0909
$ws.Range("L136").Value = 7828.6362
$ws.Range("M136").Value = -8892.090899999999
$ws.Range("N136").Value = -12928.6362

$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 3000
$ws.Range("J10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3338

# Row 136
$ws.Range("H136").Value = 2467.1912
$ws.Range("I136").Value = 2672.138
$ws.Range("J136").Value = 1278.5
$ws.Range("K136").Value = 8016.414
$ws.Range("L136").Value = 3835.5
$ws.Range("M136").Value = -5466.414
$ws.Range("N136").Value = -8935.5
